$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This block of rows (149-170 on sheet "Artfynd", excluding the untouched rows
# 155 and 158) has been re-sorted: every row's entire contents move to a
# different row position in a single 20-row cycle. Row 155 and 158 keep their
# original contents. We replicate this by staging a snapshot of each source
# row in a scratch area far below the used range, then clearing and
# value-pasting each snapshot into its final destination row. Using
# Copy/PasteSpecial(xlPasteValues) (rather than reading/writing .Value2)
# preserves the original cell typing (e.g. numeric-looking or date-looking
# text stays text instead of being reinterpreted as a number/date), and
# ClearContents before each paste removes any cell that must become blank in
# its new position (matching the source row's blank cells).

$rows = 149,150,151,152,153,154,156,157,159,160,161,162,163,164,165,166,167,168,169,170

$permMap = @{
  149 = 167
  150 = 149
  151 = 160
  152 = 156
  153 = 154
  154 = 168
  156 = 162
  157 = 166
  159 = 161
  160 = 152
  161 = 170
  162 = 163
  163 = 169
  164 = 150
  165 = 159
  166 = 151
  167 = 157
  168 = 164
  169 = 165
  170 = 153
}

$xlPasteValues = -4163
$stagingStart = 20000

# Step 1: snapshot each involved row's current (pre-edit) contents into a
# scratch row so later overwrites don't clobber data still needed as a
# source for another destination row.
$stagingMap = @{}
$i = 0
foreach ($r in $rows) {
  $stageRow = $stagingStart + $i
  $stagingMap[$r] = $stageRow
  $srcRange = $ws.Range("A" + $r + ":AY" + $r)
  $dstRange = $ws.Range("A" + $stageRow + ":AY" + $stageRow)
  $srcRange.Copy()
  $dstRange.PasteSpecial($xlPasteValues)
  $i++
}

# Step 2: for each destination row, clear it then paste in the staged
# snapshot of its source row.
foreach ($dest in $rows) {
  $srcRowOriginal = $permMap[$dest]
  $stageRow = $stagingMap[$srcRowOriginal]
  $destRange = $ws.Range("A" + $dest + ":AY" + $dest)
  $stageRange = $ws.Range("A" + $stageRow + ":AY" + $stageRow)
  $destRange.ClearContents()
  $stageRange.Copy()
  $destRange.PasteSpecial($xlPasteValues)
}

# Step 3: clean up the scratch rows.
foreach ($r in $rows) {
  $stageRow = $stagingMap[$r]
  $ws.Range("A" + $stageRow + ":AY" + $stageRow).ClearContents()
}

$excel.CutCopyMode = 0
